# NetLiquidity / FRED_Data / WTREGEN.xlsx refresh
# - Append the newest weekly observation to the "Data" sheet.
# - Update the refreshed metadata on the "SeriesInfo" sheet.

$wb = $excel.ActiveWorkbook

# ---- "Data" sheet: append new observation row 94 ----
$data = $wb.Worksheets.Item("Data")

# Duplicate the formatting of the last existing data row (row 93) onto the
# new row 94 before writing values, so the new date cell keeps the same
# date/time number format + style as the rest of column A.
$data.Range("A93:B93").Copy()
$data.Range("A94:B94").PasteSpecial(-4122)

$data.Range("A94").Value = 45119
$data.Range("B94").Value = 514.337

# ---- "SeriesInfo" sheet: refresh FRED metadata ----
$info = $wb.Worksheets.Item("SeriesInfo")

# B3/B4/B7 hold plain "YYYY-MM-DD" strings. A bare .Value assignment gets
# auto-recognised as a real date by Excel's input parser (which would change
# the cell's stored type/format), so force text entry via a "@" number
# format, then clear the format back to the sheet's default (these cells
# carry no explicit style in the source file) once the literal text value is
# committed.
$info.Range("B3").NumberFormat = "@"
$info.Range("B3").Value = "2023-07-20"
$info.Range("B3").ClearFormats()

$info.Range("B4").NumberFormat = "@"
$info.Range("B4").Value = "2023-07-20"
$info.Range("B4").ClearFormats()

$info.Range("B7").NumberFormat = "@"
$info.Range("B7").Value = "2023-07-12"
$info.Range("B7").ClearFormats()

# B14 ("2023-07-13 15:35:18-05") isn't recognised as a date/time by Excel's
# parser (the "-05" UTC-offset suffix isn't a valid literal), so it is
# stored as plain text with no special handling needed.
$info.Range("B14").Value = "2023-07-13 15:35:18-05"
